$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 54815
$ws.Range("F4").Value = 1346
$ws.Range("F6").Value = 337
$ws.Range("F8").Value = 769
$ws.Range("F9").Value = 411
$ws.Range("F10").Value = 3077
$ws.Range("F11").Value = 916
$ws.Range("F12").Value = 5246
$ws.Range("F13").Value = 1290
$ws.Range("F14").Value = 1071
$ws.Range("F16").Value = 851
$ws.Range("F18").Value = 416
$ws.Range("F19").Value = 1300
$ws.Range("F23").Value = 373
$ws.Range("F24").Value = 35
$ws.Range("F26").Value = 26
$ws.Range("F29").Value = 5145
$ws.Range("F31").Value = 5056
$ws.Range("F32").Value = 9062
$ws.Range("F34").Value = 154
$ws.Range("F35").Value = 138
$ws.Range("F36").Value = 229
$ws.Range("F37").Value = 436
$ws.Range("F38").Value = 119
$ws.Range("F39").Value = 93
$ws.Range("F40").Value = 4220
$ws.Range("F41").Value = 250

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 97
$ws.Range("F12").Value = 1137

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1346
$ws.Range("F5").Value = 337
$ws.Range("F7").Value = 769
$ws.Range("F8").Value = 411
$ws.Range("F9").Value = 916
$ws.Range("F10").Value = 97
$ws.Range("F11").Value = 1290
$ws.Range("F14").Value = 1071
$ws.Range("F16").Value = 851
$ws.Range("F17").Value = 416
$ws.Range("F19").Value = 1300
$ws.Range("F24").Value = 373
$ws.Range("F25").Value = 35
$ws.Range("F28").Value = 5145
$ws.Range("F30").Value = 9062
$ws.Range("F33").Value = 154
$ws.Range("F34").Value = 138
$ws.Range("F35").Value = 229
$ws.Range("F36").Value = 436
$ws.Range("F39").Value = 119
$ws.Range("F40").Value = 93
$ws.Range("F41").Value = 4220
$ws.Range("F48").Value = 250
